$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Uren geschat" (estimated minutes) column D values for rows 2-12
# (order chosen to match the shared-string insertion order of the target file)
$ws.Range("D5").Value = "150 minuten"
$ws.Range("D2").Value = "60 minuten"
$ws.Range("D3").Value = "45 minuten"
$ws.Range("D4").Value = "120 minuten"
$ws.Range("D6").Value = "75 minuten"
$ws.Range("D7").Value = "30 minuten"
$ws.Range("D8").Value = "15 minuten"
$ws.Range("D10").Value = "20 minuten"
$ws.Range("D11").Value = "10 minuten"
$ws.Range("D9").Value = "15 minuten"
$ws.Range("D12").Value = "15 minuten"

# Move "x" markers in rows 7, 8, 9, 10, 12 to new columns
$ws.Range("H7").ClearContents()
$ws.Range("J7").Value = "x"

$ws.Range("H8").ClearContents()
$ws.Range("J8").Value = "x"

$ws.Range("I9").ClearContents()
$ws.Range("J9").Value = "x"

$ws.Range("I10").ClearContents()
$ws.Range("J10").Value = "x"

$ws.Range("J12").ClearContents()
$ws.Range("H12").Value = "x"

# Update the active selection to K12
$ws.Range("K12").Select()
